$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8249614238739014
$ws.Range("B1").Value = 1.147810697555542
$ws.Range("C1").Value = 1.99462366104126
$ws.Range("D1").Value = 4.536254405975342
$ws.Range("E1").Value = 2.19387412071228
